$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confusion-matrix counts (TN/FP row and FN/TP row)
$ws.Range("G3").Value = 96
$ws.Range("H3").Value = 4
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 19

# New helper cell with a small formula
$ws.Range("K8").Formula = "=1/100"

# Bump the displayed precision of the computed metrics from 2 to 3 decimals
$fmt3 = "_-* #,##0.000_-;\-* #,##0.000_-;_-* ""-""??_-;_-@_-"
$ws.Range("D8:D9").NumberFormat = $fmt3
$ws.Range("D10:D13").NumberFormat = $fmt3

# Match the reported selection after the edits
$ws.Range("D8:D9").Select()
